$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.693.10"
$ws.Range("E2").Value = "'  -0.70%  "
$ws.Range("D3").Value = "'2.528.11"
$ws.Range("E3").Value = "'  -1.83%  "
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("D5").Value = "'309.17"
$ws.Range("E5").Value = "'  -1.89%  "
$ws.Range("D6").Value = "'100.77"
$ws.Range("E6").Value = "'  +0.81%  "
$ws.Range("D7").Value = "'0.567"
$ws.Range("E7").Value = "'  -1.54%  "
$ws.Range("E8").Value = "'  +0.11%  "
$ws.Range("D9").Value = "'0.524"
$ws.Range("E9").Value = "'  -2.68%  "
$ws.Range("E10").Value = "'  -1.24%  "
$ws.Range("D11").Value = "'0.0804"
$ws.Range("E11").Value = "'  -1.11%  "
$ws.Range("D12").Value = "'7.33"
$ws.Range("E12").Value = "'  -3.19%  "
$ws.Range("E13").Value = "'  +0.22%  "
$ws.Range("D14").Value = "'2.916.41"
$ws.Range("E14").Value = "'  -1.85%  "
$ws.Range("D15").Value = "'15.40"
$ws.Range("E15").Value = "'  -1.94%  "
$ws.Range("D16").Value = "'2.516.63"
$ws.Range("E16").Value = "'  -4.89%  "
$ws.Range("E17").Value = "'  -4.18%  "
$ws.Range("D18").Value = "'42.678.42"
$ws.Range("E18").Value = "'  -0.80%  "
$ws.Range("D19").Value = "'6.72"
$ws.Range("E19").Value = "'  -2.22%  "
$ws.Range("B20").Value = "'InternetComputer(DFINITY)"
$ws.Range("C20").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "'12.32"
$ws.Range("E20").Value = "'  -2.78%  "
$ws.Range("B21").Value = "'ShibaInu"
$ws.Range("C21").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "'0.0₃0950"
$ws.Range("E21").Value = "'  -2.08%  "
$ws.Range("D22").Value = "'69.56"
$ws.Range("E22").Value = "'  +0.20%  "
$ws.Range("D23").Value = "'243.86"
$ws.Range("E23").Value = "'  -2.69%  "
$ws.Range("E24").Value = "'  -2.67%  "
$ws.Range("E25").Value = "'  -3.27%  "
$ws.Range("E26").Value = "'  -0.01%  "
$ws.Range("D27").Value = "'25.47"
$ws.Range("E27").Value = "'  -6.14%  "
$ws.Range("E28").Value = "'  -3.02%  "
$ws.Range("D29").Value = "'10.13"
$ws.Range("E29").Value = "'  -1.88%  "
$ws.Range("D30").Value = "'38.77"
$ws.Range("E30").Value = "'  -4.53%  "
$ws.Range("D31").Value = "'157.75"
$ws.Range("E31").Value = "'  +0.08%  "
$ws.Range("D32").Value = "'5.74"
$ws.Range("E32").Value = "'  -2.03%  "
$ws.Range("E33").Value = "'  +10.98%  "
$ws.Range("D34").Value = "'0.0786"
$ws.Range("E34").Value = "'  -2.44%  "
$ws.Range("E35").Value = "'  -1.53%  "
$ws.Range("D36").Value = "'3.16"
$ws.Range("E36").Value = "'  -8.39%  "
$ws.Range("E37").Value = "'  -7.05%  "
$ws.Range("D38").Value = "'17.79"
$ws.Range("E38").Value = "'  -5.61%  "
$ws.Range("E39").Value = "'  -1.50%  "
$ws.Range("E40").Value = "'  -0.90%  "
$ws.Range("D41").Value = "'4.21"
$ws.Range("E41").Value = "'  +3.95%  "
$ws.Range("D42").Value = "'21.94"
$ws.Range("E42").Value = "'  -8.08%  "
$ws.Range("E43").Value = "'  +0.13%  "
$ws.Range("E44").Value = "'  -1.55%  "
$ws.Range("E45").Value = "'  +0.67%  "
$ws.Range("D46").Value = "'2.007.98"
$ws.Range("E46").Value = "'  -0.10%  "
$ws.Range("D47").Value = "'8.89"
$ws.Range("E47").Value = "'  -0.39%  "
$ws.Range("D48").Value = "'2.771.61"
$ws.Range("E48").Value = "'  -1.82%  "
$ws.Range("E49").Value = "'  -4.02%  "
$ws.Range("B50").Value = "'ordi"
$ws.Range("C50").Value = "'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").Value = "'72.24"
$ws.Range("E50").Value = "'  -3.74%  "
$ws.Range("B51").Value = "'BitcoinSV"
$ws.Range("C51").Value = "'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").Value = "'79.11"
$ws.Range("E51").Value = "'  -3.82%  "
